$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column F (dSF), rows 2-32
$values = @{
    2  = -6
    3  = -6
    4  = 1
    5  = -1
    6  = 5
    7  = 1
    8  = 4
    9  = -5
    10 = 10
    11 = 1
    12 = 2
    13 = 2
    14 = 2
    15 = -1
    16 = -2
    17 = 2
    18 = -2
    19 = 6
    20 = 7
    21 = 7
    22 = -1
    23 = -3
    24 = 3
    25 = 5
    26 = -2
    27 = -4
    28 = 2
    29 = 12
    30 = 4
    31 = -4
    32 = 5
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
